$d = $word.ActiveDocument
$d.Content.Find.Execute("2.a, 3.b, 3.g i 4.f", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.b, 3.g i 4.f", 2)
